# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Sat Dec 23 07:09:53 UTC 2023 with GitHub Actions".
#
# Column D ("Price") holds numbers formatted with a dotted thousands
# separator (e.g. "43.646.22") alongside plain decimals (e.g. "96.88").
# The sheet stores every one of these as literal text, so any Price cell
# whose new value would otherwise be auto-parsed as a number by Excel is
# first stamped with a Text number format ("@") to keep it text (this also
# preserves trailing zeros such as "0.850").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2"; Value = "43.646.22"; AsText = $false },
    @{ Addr = "E2"; Value = "  -0.59%  "; AsText = $false },
    @{ Addr = "D3"; Value = "2.293.40"; AsText = $false },
    @{ Addr = "E3"; Value = "  +1.26%  "; AsText = $false },
    @{ Addr = "E4"; Value = "  +0.07%  "; AsText = $false },
    @{ Addr = "D5"; Value = "96.88"; AsText = $true },
    @{ Addr = "E5"; Value = "  +2.16%  "; AsText = $false },
    @{ Addr = "D6"; Value = "266.71"; AsText = $true },
    @{ Addr = "E6"; Value = "  -1.50%  "; AsText = $false },
    @{ Addr = "D7"; Value = "0.624"; AsText = $true },
    @{ Addr = "E7"; Value = "  -0.28%  "; AsText = $false },
    @{ Addr = "E8"; Value = "  +0.01%  "; AsText = $false },
    @{ Addr = "D9"; Value = "0.612"; AsText = $true },
    @{ Addr = "E9"; Value = "  -2.13%  "; AsText = $false },
    @{ Addr = "D10"; Value = "45.79"; AsText = $true },
    @{ Addr = "E10"; Value = "  -2.53%  "; AsText = $false },
    @{ Addr = "D11"; Value = "0.0936"; AsText = $true },
    @{ Addr = "E12"; Value = "  -3.07%  "; AsText = $false },
    @{ Addr = "E13"; Value = "  +0.34%  "; AsText = $false },
    @{ Addr = "D14"; Value = "2.636.64"; AsText = $false },
    @{ Addr = "E14"; Value = "  +1.28%  "; AsText = $false },
    @{ Addr = "D15"; Value = "15.16"; AsText = $true },
    @{ Addr = "E15"; Value = "  -1.19%  "; AsText = $false },
    @{ Addr = "D16"; Value = "0.850"; AsText = $true },
    @{ Addr = "E16"; Value = "  +2.87%  "; AsText = $false },
    @{ Addr = "D17"; Value = "2.298.10"; AsText = $false },
    @{ Addr = "E17"; Value = "  +2.44%  "; AsText = $false },
    @{ Addr = "D18"; Value = "43.606.54"; AsText = $false },
    @{ Addr = "E18"; Value = "  -0.57%  "; AsText = $false },
    @{ Addr = "E19"; Value = "  +2.06%  "; AsText = $false },
    @{ Addr = "D20"; Value = "6.17"; AsText = $true },
    @{ Addr = "E20"; Value = "  +0.16%  "; AsText = $false },
    @{ Addr = "D21"; Value = "71.81"; AsText = $true },
    @{ Addr = "E21"; Value = "  +1.19%  "; AsText = $false },
    @{ Addr = "D22"; Value = "2.41"; AsText = $true },
    @{ Addr = "E22"; Value = "  +5.14%  "; AsText = $false },
    @{ Addr = "D23"; Value = "232.79"; AsText = $true },
    @{ Addr = "E23"; Value = "  -1.19%  "; AsText = $false },
    @{ Addr = "D24"; Value = "9.16"; AsText = $true },
    @{ Addr = "E24"; Value = "  -8.82%  "; AsText = $false },
    @{ Addr = "E25"; Value = "  +0.04%  "; AsText = $false },
    @{ Addr = "D26"; Value = "2.49"; AsText = $true },
    @{ Addr = "E26"; Value = "  -0.34%  "; AsText = $false },
    @{ Addr = "D27"; Value = "11.17"; AsText = $true },
    @{ Addr = "E27"; Value = "  -1.45%  "; AsText = $false },
    @{ Addr = "D28"; Value = "3.48"; AsText = $true },
    @{ Addr = "E28"; Value = "  +3.39%  "; AsText = $false },
    @{ Addr = "D29"; Value = "40.02"; AsText = $true },
    @{ Addr = "E29"; Value = "  +0.71%  "; AsText = $false },
    @{ Addr = "D30"; Value = "2.22"; AsText = $true },
    @{ Addr = "E30"; Value = "  -1.80%  "; AsText = $false },
    @{ Addr = "B31"; Value = "EthereumClassic"; AsText = $false },
    @{ Addr = "C31"; Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; AsText = $false },
    @{ Addr = "D31"; Value = "22.31"; AsText = $true },
    @{ Addr = "E31"; Value = "  +1.17%  "; AsText = $false },
    @{ Addr = "B32"; Value = "Monero"; AsText = $false },
    @{ Addr = "C32"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; AsText = $false },
    @{ Addr = "D32"; Value = "175.49"; AsText = $true },
    @{ Addr = "E32"; Value = "  +1.25%  "; AsText = $false },
    @{ Addr = "D33"; Value = "0.0883"; AsText = $true },
    @{ Addr = "E33"; Value = "  -2.72%  "; AsText = $false },
    @{ Addr = "D34"; Value = "5.37"; AsText = $true },
    @{ Addr = "E34"; Value = "  -3.93%  "; AsText = $false },
    @{ Addr = "E35"; Value = "  +1.11%  "; AsText = $false },
    @{ Addr = "E36"; Value = "  -2.69%  "; AsText = $false },
    @{ Addr = "E37"; Value = "  +1.45%  "; AsText = $false },
    @{ Addr = "D38"; Value = "4.37"; AsText = $true },
    @{ Addr = "E38"; Value = "  -0.89%  "; AsText = $false },
    @{ Addr = "D39"; Value = "3.43"; AsText = $true },
    @{ Addr = "E39"; Value = "  -0.84%  "; AsText = $false },
    @{ Addr = "D40"; Value = "0.239"; AsText = $true },
    @{ Addr = "E40"; Value = "  -4.43%  "; AsText = $false },
    @{ Addr = "D41"; Value = "2.33"; AsText = $true },
    @{ Addr = "E41"; Value = "  +5.53%  "; AsText = $false },
    @{ Addr = "D42"; Value = "12.31"; AsText = $true },
    @{ Addr = "E42"; Value = "  -0.68%  "; AsText = $false },
    @{ Addr = "D43"; Value = "1.35"; AsText = $true },
    @{ Addr = "E43"; Value = "  +12.83%  "; AsText = $false },
    @{ Addr = "D44"; Value = "64.59"; AsText = $true },
    @{ Addr = "E44"; Value = "  +5.38%  "; AsText = $false },
    @{ Addr = "E45"; Value = "  +3.53%  "; AsText = $false },
    @{ Addr = "E46"; Value = "  -4.06%  "; AsText = $false },
    @{ Addr = "E47"; Value = "  -0.28%  "; AsText = $false },
    @{ Addr = "D48"; Value = "98.19"; AsText = $true },
    @{ Addr = "E48"; Value = "  -1.42%  "; AsText = $false },
    @{ Addr = "D49"; Value = "1.19"; AsText = $true },
    @{ Addr = "E49"; Value = "  +0.63%  "; AsText = $false },
    @{ Addr = "D50"; Value = "2.516.11"; AsText = $false },
    @{ Addr = "E50"; Value = "  +1.28%  "; AsText = $false },
    @{ Addr = "E51"; Value = "  +0.08%  "; AsText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    if ($u.AsText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}

